$wb = $excel.ActiveWorkbook

# --- responseTime sheet ---
$ws = $wb.Worksheets.Item("responseTime")
$ws.Range("C2").Value = 9.27
$ws.Range("B3").Value = 2.81
$ws.Range("C3").Value = 12.62
$ws.Range("B4").Value = 3.63
$ws.Range("C4").Value = 18.09
$ws.Range("B5").Value = 4.44
$ws.Range("C5").Value = 24.72
$ws.Range("B6").Value = 6.03
$ws.Range("C6").Value = 29.55

# --- requestsPerSecond sheet ---
$ws = $wb.Worksheets.Item("requestsPerSecond")
$ws.Range("B2").Value = 395.03
$ws.Range("C2").Value = 107.7
$ws.Range("B3").Value = 707.25
$ws.Range("C3").Value = 158.07
$ws.Range("B4").Value = 821.93
$ws.Range("C4").Value = 165.35
$ws.Range("B5").Value = 894.12
$ws.Range("C5").Value = 160.93
$ws.Range("B6").Value = 824.47
$ws.Range("C6").Value = 168.45

# --- cpuUsage sheet ---
$ws = $wb.Worksheets.Item("cpuUsage")
$ws.Range("B2").Value = 71.09
$ws.Range("C2").Value = 79.95
$ws.Range("B3").Value = 62.85
$ws.Range("C3").Value = 94.57
$ws.Range("B4").Value = 71.81
$ws.Range("C4").Value = 97.93
$ws.Range("B5").Value = 69.87
$ws.Range("C5").Value = 96.96
$ws.Range("B6").Value = 67.42
$ws.Range("C6").Value = 97.76

# --- memoryUsage sheet ---
$ws = $wb.Worksheets.Item("memoryUsage")
$ws.Range("B2").Value = 271.44
$ws.Range("C2").Value = 230.69
$ws.Range("B3").Value = 294.48
$ws.Range("C3").Value = 295.98
$ws.Range("B4").Value = 224.71
$ws.Range("C4").Value = 258.43
$ws.Range("B5").Value = 225.65
$ws.Range("C5").Value = 368.49
$ws.Range("B6").Value = 235.35
$ws.Range("C6").Value = 263.88

# --- Selection on memoryUsage stays put, but active tab moves to Charts ---
$wsCharts = $wb.Worksheets.Item("Charts")
$wsCharts.Activate()
$wsCharts.Range("S23").Select()
